{"js": "// Adds two new paragraphs to the end of the document body, right after the\n// \"Belgenin 3. paragraf\u0131\" paragraph:\n//   1) \"Belgemin 4. paragraf\u0131\"\n//   2) Five runs reading (concatenated):\n//        \"-baba ve o\u011fulun hikayesiBurada bize verilen g\u00f6rev hikayen\u0131n bu k\u0131sm\u0131n\u0131 tamamlamak \"\n//\n// All new runs/paragraphs reuse the Times New Roman / 24pt (sz 48) formatting\n// already used throughout the document, and the paragraphs use the \"Normal\"\n// style - exactly like the author's OOXML diff.\n//\n// We build the new paragraphs via a literal OOXML (Flat OPC) fragment and\n// insert it with Range.insertOoxml so that the five runs of the second\n// paragraph remain distinct <w:r> elements (as in the diff) instead of being\n// silently merged together by identical-formatting run coalescing.\n\nconst rPr =\n  '<w:rPr>' +\n  '<w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:eastAsia=\"Times New Roman\" w:cs=\"Times New Roman\"/>' +\n  '<w:sz w:val=\"48\"/>' +\n  '<w:szCs w:val=\"48\"/>' +\n  '</w:rPr>';\n\nfunction run(text) {\n  const preserve = /^\\s|\\s$/.test(text) ? ' xml:space=\"preserve\"' : '';\n  return '<w:r>' + rPr + '<w:t' + preserve + '>' + text + '</w:t></w:r>';\n}\n\nconst pPr =\n  '<w:pPr>' +\n  '<w:pStyle w:val=\"Normal\"/>' +\n  rPr +\n  '</w:pPr>';\n\nconst paragraph4 = '<w:p>' + pPr + run('Belgemin 4. paragraf\u0131') + '</w:p>';\n\nconst paragraph5Runs = [\n  '-baba ve o\u011fulun ',\n  'hikayesiBurada',\n  ' bize verilen g\u00f6rev ',\n  'hikayen\u0131n',\n  ' bu k\u0131sm\u0131n\u0131 tamamlamak ',\n]\n  .map(run)\n  .join('');\n\nconst paragraph5 = '<w:p>' + pPr + paragraph5Runs + '</w:p>';\n\nconst flatOpc =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  paragraph4 +\n  paragraph5 +\n  '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '</pkg:package>';\n\nconst body = context.document.body;\nbody.paragraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = body.paragraphs.items[body.paragraphs.items.length - 1];\nconst insertionRange = lastParagraph.getRange();\ninsertionRange.insertOoxml(flatOpc, Word.InsertLocation.after);\n\nawait context.sync();\n", "ps1": "# Adds two new paragraphs to the end of the document, right after the\n# \"Belgenin 3. paragraf\u0131\" paragraph:\n#   1) \"Belgemin 4. paragraf\u0131\"\n#   2) Five runs reading (concatenated):\n#        \"-baba ve o\u011fulun hikayesiBurada bize verilen g\u00f6rev hikayen\u0131n bu k\u0131sm\u0131n\u0131 tamamlamak \"\n#\n# All new runs/paragraphs reuse the Times New Roman / 24pt (sz 48) formatting\n# already used throughout the document, and the paragraphs use the \"Normal\"\n# style - exactly like the author's OOXML diff.\n#\n# The new paragraph content is built as a literal WordprocessingML fragment\n# and inserted with Range.InsertXML so the five runs of the second paragraph\n# remain distinct <w:r> elements (as in the diff) instead of being silently\n# coalesced together by identical-formatting run merging.\n\n$d = $word.ActiveDocument\n\n$lastParagraph = $d.Paragraphs.Item($d.Paragraphs.Count)\n$insertionRange = $lastParagraph.Range\n\n$wordMlNs = 'xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"'\n\n$runProps = '<w:rPr><w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:eastAsia=\"Times New Roman\" w:cs=\"Times New Roman\"/><w:sz w:val=\"48\"/><w:szCs w:val=\"48\"/></w:rPr>'\n$paraProps = '<w:pPr><w:pStyle w:val=\"Normal\"/>' + $runProps + '</w:pPr>'\n\nfunction New-Run($text) {\n    $preserve = \"\"\n    if ($text -match '^\\s' -or $text -match '\\s$') {\n        $preserve = ' xml:space=\"preserve\"'\n    }\n    return '<w:r>' + $runProps + '<w:t' + $preserve + '>' + $text + '</w:t></w:r>'\n}\n\n$paragraph4 = '<w:p ' + $wordMlNs + '>' + $paraProps + (New-Run \"Belgemin 4. paragraf\u0131\") + '</w:p>'\n\n$paragraph5RunTexts = @(\n    \"-baba ve o\u011fulun \",\n    \"hikayesiBurada\",\n    \" bize verilen g\u00f6rev \",\n    \"hikayen\u0131n\",\n    \" bu k\u0131sm\u0131n\u0131 tamamlamak \"\n)\n\n$paragraph5Runs = \"\"\nforeach ($runText in $paragraph5RunTexts) {\n    $paragraph5Runs += New-Run $runText\n}\n$paragraph5 = '<w:p ' + $wordMlNs + '>' + $paraProps + $paragraph5Runs + '</w:p>'\n\n$newParagraphsXml = $paragraph4 + $paragraph5\n\n$insertionRange.InsertXML($newParagraphsXml, \"After\")\n"}
